$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2166.1538
$ws.Range("I17").Value = 796
$ws.Range("K17").Value = 2388
$ws.Range("M17").Value = -2220
$ws.Range("H99").Value = 128.33333
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H116").Value = 4319.125
$ws.Range("I116").Value = 3691
$ws.Range("J116").Value = 4696
$ws.Range("K116").Value = 3691
$ws.Range("L116").Value = 4696
$ws.Range("M116").Value = -249
$ws.Range("N116").Value = -11580
$ws.Range("H132").Value = 14609.6875
$ws.Range("I132").Value = 14583.667
$ws.Range("K132").Value = 43751.001
$ws.Range("M132").Value = -41221.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 100
$ws.Range("I5").Value = 100
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 100
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 12
$ws.Range("N5").ClearContents()
$ws.Range("H6").Value = 9002450
$ws.Range("I6").Value = 7503000
$ws.Range("J6").Value = 15000250
$ws.Range("K6").Value = 7503000
$ws.Range("L6").Value = 15000250
$ws.Range("M6").Value = -7502827
$ws.Range("N6").Value = -15000596
$ws.Range("H25").Value = 887.25
$ws.Range("I25").Value = 887.25
$ws.Range("K25").Value = 887.25
$ws.Range("M25").Value = -485.25
$ws.Range("H45").Value = 3121.5
$ws.Range("I45").Value = 2055.6
$ws.Range("K45").Value = 2055.6
$ws.Range("M45").Value = -1678.6
$ws.Range("H61").Value = 2714.1428
$ws.Range("I61").Value = 2714.1428
$ws.Range("K61").Value = 2714.1428
$ws.Range("M61").Value = -2502.1428
$ws.Range("H97").Value = 752.1818
$ws.Range("I97").Value = 844.44446
$ws.Range("K97").Value = 844.44446
$ws.Range("M97").Value = -348.44446
$ws.Range("H136").Value = 2714.1428
$ws.Range("I136").Value = 2714.1428
$ws.Range("K136").Value = 8142.428400000001
$ws.Range("M136").Value = -5592.428400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 100
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 15
$ws.Range("N4").ClearContents()
$ws.Range("H80").Value = 1011.6
$ws.Range("I80").Value = 646
$ws.Range("J80").Value = 1377.2
$ws.Range("K80").Value = 646
$ws.Range("L80").Value = 1377.2
$ws.Range("M80").Value = 352
$ws.Range("N80").Value = -3373.2
$ws.Range("H83").Value = 1011.6
$ws.Range("I83").Value = 646
$ws.Range("J83").Value = 1377.2
$ws.Range("K83").Value = 3230
$ws.Range("L83").Value = 6886
$ws.Range("M83").Value = 1762
$ws.Range("N83").Value = -16870

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 68.36364
$ws.Range("J7").Value = 81.166664
$ws.Range("L7").Value = 81.166664
$ws.Range("N7").Value = -307.166664
$ws.Range("H22").Value = 4444
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 4444
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 4444
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -5144
$ws.Range("H31").Value = 4999.533
$ws.Range("I31").Value = 2165.5
$ws.Range("K31").Value = 2165.5
$ws.Range("M31").Value = -1870.5
$ws.Range("H34").Value = 4999.533
$ws.Range("I34").Value = 2165.5
$ws.Range("K34").Value = 2165.5
$ws.Range("M34").Value = -1963.5
$ws.Range("H96").Value = 15333.333
$ws.Range("J96").Value = 15333.333
$ws.Range("L96").Value = 15333.333
$ws.Range("N96").Value = -20825.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 787.6429000000001
$ws.Range("I34").Value = 170.4
$ws.Range("J34").Value = 1130.5555
$ws.Range("K34").Value = 511.2
$ws.Range("L34").Value = 3391.6665
$ws.Range("M34").Value = -427.2
$ws.Range("N34").Value = -3559.6665
$ws.Range("H107").Value = 476.5484
$ws.Range("I107").Value = 218.2
$ws.Range("K107").Value = 654.5999999999999
$ws.Range("M107").Value = 1265.4
$ws.Range("H113").Value = 1304.4
$ws.Range("I113").Value = 900
$ws.Range("J113").Value = 1349.3334
$ws.Range("K113").Value = 2700
$ws.Range("L113").Value = 4048.0002
$ws.Range("M113").Value = -530
$ws.Range("N113").Value = -8388.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 156.51428
$ws.Range("I2").Value = 43.517242
$ws.Range("J2").Value = 702.6667
$ws.Range("K2").Value = 43.517242
$ws.Range("L2").Value = 702.6667
$ws.Range("M2").Value = 69.48275799999999
$ws.Range("N2").Value = -928.6667
$ws.Range("H7").Value = 28000376
$ws.Range("I7").Value = 17778278
$ws.Range("K7").Value = 17778278
$ws.Range("M7").Value = -17778166
$ws.Range("H8").Value = 28000376
$ws.Range("I8").Value = 17778278
$ws.Range("K8").Value = 17778278
$ws.Range("M8").Value = -17778139
$ws.Range("H14").Value = 503
$ws.Range("I14").Value = 503
$ws.Range("K14").Value = 503
$ws.Range("M14").Value = -335
$ws.Range("H80").Value = 2615.7144
$ws.Range("I80").Value = 2385
$ws.Range("K80").Value = 2385
$ws.Range("M80").Value = -1387
$ws.Range("H83").Value = 2615.7144
$ws.Range("I83").Value = 2385
$ws.Range("K83").Value = 11925
$ws.Range("M83").Value = -6933
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H126").Value = 4080
$ws.Range("I126").Value = 4116.2
$ws.Range("K126").Value = 12348.6
$ws.Range("M126").Value = -9878.599999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 863.6667
$ws.Range("I22").Value = 697.25
$ws.Range("J22").Value = 996.8
$ws.Range("K22").Value = 697.25
$ws.Range("L22").Value = 996.8
$ws.Range("M22").Value = -402.25
$ws.Range("N22").Value = -1586.8
$ws.Range("H27").Value = 863.6667
$ws.Range("I27").Value = 697.25
$ws.Range("J27").Value = 996.8
$ws.Range("K27").Value = 697.25
$ws.Range("L27").Value = 996.8
$ws.Range("M27").Value = -590.25
$ws.Range("N27").Value = -1210.8
$ws.Range("H46").Value = 6346.067
$ws.Range("I46").Value = 2097.5
$ws.Range("J46").Value = 7891
$ws.Range("K46").Value = 2097.5
$ws.Range("L46").Value = 7891
$ws.Range("M46").Value = -1909.5
$ws.Range("N46").Value = -8267
$ws.Range("H81").Value = 19000
$ws.Range("J81").Value = 19000
$ws.Range("L81").Value = 19000
$ws.Range("N81").Value = -20996
$ws.Range("H82").Value = 2819.8
$ws.Range("I82").Value = 516.6667
$ws.Range("K82").Value = 516.6667
$ws.Range("M82").Value = -155.6667
$ws.Range("H84").Value = 19000
$ws.Range("J84").Value = 19000
$ws.Range("L84").Value = 57000
$ws.Range("N84").Value = -66984
$ws.Range("H85").Value = 2819.8
$ws.Range("I85").Value = 516.6667
$ws.Range("K85").Value = 516.6667
$ws.Range("M85").Value = 731.3333
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()
$ws.Range("H132").Value = 3724.25
$ws.Range("I132").Value = 3898.5
$ws.Range("K132").Value = 11695.5
$ws.Range("M132").Value = -9165.5
$ws.Range("H136").Value = 5000.6665
$ws.Range("I136").Value = 5000.6665
$ws.Range("K136").Value = 15001.9995
$ws.Range("M136").Value = -12451.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 233.33333
$ws.Range("I2").Value = 200
$ws.Range("J2").Value = 400
$ws.Range("K2").Value = 200
$ws.Range("L2").Value = 400
$ws.Range("M2").Value = -88
$ws.Range("N2").Value = -624
$ws.Range("H80").Value = 59518.8
$ws.Range("J80").Value = 59518.8
$ws.Range("L80").Value = 59518.8
$ws.Range("N80").Value = -61514.8
$ws.Range("H83").Value = 59518.8
$ws.Range("J83").Value = 59518.8
$ws.Range("L83").Value = 178556.4
$ws.Range("N83").Value = -188540.4
$ws.Range("H122").Value = 2378.4092
$ws.Range("I122").Value = 2141
$ws.Range("K122").Value = 6423
$ws.Range("M122").Value = -3973
$ws.Range("H126").Value = 3659.6
$ws.Range("I126").Value = 1489.4
$ws.Range("K126").Value = 4468.200000000001
$ws.Range("M126").Value = -1998.200000000001
